$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "NTN" in G1 (merged NTE + NTW)
$ws.Range("G1").Value = "NTN"

# Fill G2:G18 with the sum of NTE (column C) and NTW (column D) for each row
for ($r = 2; $r -le 18; $r++) {
    $nte = [double]($ws.Cells.Item($r, 3).Value2)
    $ntw = [double]($ws.Cells.Item($r, 4).Value2)
    $ws.Cells.Item($r, 7).Value2 = $nte + $ntw
}

# Update the selected cell to match the new active selection
$ws.Range("J14").Select()
